# Commit: "Change Excel Field View to Cache, And set default value to FALSE"
#
# The "Property" sheet has a header row (row 1) naming each column, and a
# column "View" (F) whose default values across all data rows were TRUE.
# This edit renames the header "View" -> "Cache", and resets every data
# row's value in that column from TRUE to FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header of column F from "View" to "Cache".
$ws.Range("F1").Value = "Cache"

# Set the default value of column F (rows 2-25) to FALSE.
$ws.Range("F2:F25").Value = $false

# Reflect the edited range as the active selection (F2 -> F2:F25), matching
# the state Excel leaves behind after editing that column.
[void]$ws.Range("F2:F25").Select()
